# Apply renamed sheets, rerun pooled models with normalized d2c (DistCenter -> DistCenter_pc)
$wb = $excel.ActiveWorkbook

# --- Step 1: rename sheets to temporary unique names to avoid collisions ---
for ($i = 1; $i -le 9; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = "__tmp" + $i + "__"
}

# --- Step 2: rename sheets to their final names (in workbook/tab order) ---
$finalNames = @("summ2","summ5","summ3","summ0","summ14","summ1","summ4","summ7","summ11")
for ($i = 1; $i -le 9; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $finalNames[$i-1]
}

# --- Step 3: update the regression results (coefficient / p columns) per sheet ---

$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = -0.918999159934639
$ws.Range("C2").Value = [double]"3.467030490610146e-07"
$ws.Range("B3").Value = -0.6033865476785383
$ws.Range("C3").Value = [double]"3.219578827003218e-38"
$ws.Range("B4").Value = 0.2738140163022489
$ws.Range("C4").Value = [double]"1.590199377767316e-77"
$ws.Range("B5").Value = 0.0005550828642443791
$ws.Range("C5").Value = 0
$ws.Range("B6").Value = 0.01386435498940276
$ws.Range("C6").Value = [double]"1.715239838836841e-32"
$ws.Range("B7").Value = 0.009005470086171615
$ws.Range("C7").Value = 0.7803609577960849
$ws.Range("B8").Value = 0.4475724572348374
$ws.Range("C8").Value = [double]"3.956609165367508e-21"
$ws.Range("B9").Value = 0.2919163640731889
$ws.Range("C9").Value = [double]"2.303554398906974e-06"
$ws.Range("B10").Value = [double]"-7.388605036827405e-05"
$ws.Range("C10").Value = [double]"7.678372064705444e-123"
$ws.Range("B11").Value = -0.01263658547111004
$ws.Range("C11").Value = 0.179205660948899
$ws.Range("A12").Value = "DistCenter_pc"
$ws.Range("B12").Value = 0.002923176953791801
$ws.Range("C12").Value = [double]"1.130049007761167e-13"
$ws.Range("B13").Value = 0.3608054233951324
$ws.Range("C13").Value = 0.0250341209765402
$ws.Range("B14").Value = 0.0002855739853420207
$ws.Range("C14").Value = 0.8259856284524387
$ws.Range("B15").Value = -0.002751817709267548
$ws.Range("C15").Value = 0.002966542320272509
$ws.Range("B16").Value = 0.6968912964808137
$ws.Range("C16").Value = [double]"4.211046784362147e-10"
$ws.Range("B17").Value = -0.6982253027379975
$ws.Range("C17").Value = [double]"1.440743358917774e-06"
$ws.Range("B18").Value = -0.0001214450241088899
$ws.Range("C18").Value = [double]"2.45551165690393e-24"

$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = -0.837617338224926
$ws.Range("C2").Value = [double]"3.615707302873595e-06"
$ws.Range("B3").Value = -0.611058413498311
$ws.Range("C3").Value = [double]"2.34636671918417e-39"
$ws.Range("B4").Value = 0.2722684947427278
$ws.Range("C4").Value = [double]"3.31695878060707e-76"
$ws.Range("B5").Value = 0.0005587833371663015
$ws.Range("C5").Value = 0
$ws.Range("B6").Value = 0.01424753398559522
$ws.Range("C6").Value = [double]"2.943560178067979e-34"
$ws.Range("B7").Value = 0.006038888346354018
$ws.Range("C7").Value = 0.8519120162813656
$ws.Range("B8").Value = 0.4445102255348063
$ws.Range("C8").Value = [double]"4.998527594035003e-21"
$ws.Range("B9").Value = 0.2691466994409506
$ws.Range("C9").Value = [double]"1.218072574974409e-05"
$ws.Range("B10").Value = [double]"-7.182037502243807e-05"
$ws.Range("C10").Value = [double]"4.996248766852013e-117"
$ws.Range("B11").Value = -0.01752797909861065
$ws.Range("C11").Value = 0.06127067369075299
$ws.Range("A12").Value = "DistCenter_pc"
$ws.Range("B12").Value = 0.003139176823171713
$ws.Range("C12").Value = [double]"1.741223540817838e-15"
$ws.Range("B13").Value = 0.3826342143493172
$ws.Range("C13").Value = 0.01758870333264658
$ws.Range("B14").Value = -0.0002268717019910899
$ws.Range("C14").Value = 0.8614532158547316
$ws.Range("B15").Value = -0.00353677922781616
$ws.Range("C15").Value = 0.0001461778679695982
$ws.Range("B16").Value = 0.6829990125643406
$ws.Range("C16").Value = [double]"8.680891679157355e-10"
$ws.Range("B17").Value = -0.6928338298643026
$ws.Range("C17").Value = [double]"1.666570542899262e-06"
$ws.Range("B18").Value = -0.00012130071689933
$ws.Range("C18").Value = [double]"1.526224422646137e-24"

$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = -0.7894135675120735
$ws.Range("C2").Value = [double]"1.226303129966875e-05"
$ws.Range("B3").Value = -0.6083459829862496
$ws.Range("C3").Value = [double]"7.534186862319907e-39"
$ws.Range("B4").Value = 0.2621947778310082
$ws.Range("C4").Value = [double]"1.554772839863902e-70"
$ws.Range("B5").Value = 0.0005550559109987965
$ws.Range("C5").Value = 0
$ws.Range("B6").Value = 0.01370337354242208
$ws.Range("C6").Value = [double]"8.61833715057318e-32"
$ws.Range("B7").Value = 0.0282600651470394
$ws.Range("C7").Value = 0.3805590286858336
$ws.Range("B8").Value = 0.421757262518809
$ws.Range("C8").Value = [double]"5.139633367548454e-19"
$ws.Range("B9").Value = 0.2502392483217661
$ws.Range("C9").Value = [double]"4.694980231305827e-05"
$ws.Range("B10").Value = [double]"-7.47378767072789e-05"
$ws.Range("C10").Value = [double]"2.249142366656074e-126"
$ws.Range("B11").Value = -0.01855241632509
$ws.Range("C11").Value = 0.04746767010511512
$ws.Range("A12").Value = "DistCenter_pc"
$ws.Range("B12").Value = 0.003039899154072705
$ws.Range("C12").Value = [double]"1.06129319162436e-14"
$ws.Range("B13").Value = 0.2443229333261155
$ws.Range("C13").Value = 0.1284584807416608
$ws.Range("B14").Value = 0.0004411498214823105
$ws.Range("C14").Value = 0.7335518072400298
$ws.Range("B15").Value = -0.003042819628520567
$ws.Range("C15").Value = 0.001101994099780205
$ws.Range("B16").Value = 0.6198244452414531
$ws.Range("C16").Value = [double]"2.512938616691911e-08"
$ws.Range("B17").Value = -0.7754940171949163
$ws.Range("C17").Value = [double]"7.622434552164802e-08"
$ws.Range("B18").Value = -0.0001121687615573495
$ws.Range("C18").Value = [double]"3.121743941084366e-21"

$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = -1.058210621547087
$ws.Range("C2").Value = [double]"4.269454761387866e-09"
$ws.Range("B3").Value = -0.5982576794899733
$ws.Range("C3").Value = [double]"5.868709551871192e-38"
$ws.Range("B4").Value = 0.2618964512557267
$ws.Range("C4").Value = [double]"1.327922054937891e-71"
$ws.Range("B5").Value = 0.0005543082273815208
$ws.Range("C5").Value = 0
$ws.Range("B6").Value = 0.01488547853727818
$ws.Range("C6").Value = [double]"2.342221879477096e-37"
$ws.Range("B7").Value = 0.02690095012268427
$ws.Range("C7").Value = 0.4042189415694231
$ws.Range("B8").Value = 0.4306835253381808
$ws.Range("C8").Value = [double]"7.316782011111174e-20"
$ws.Range("B9").Value = 0.2364469947448964
$ws.Range("C9").Value = 0.0001203378144384555
$ws.Range("B10").Value = [double]"-7.375396913576599e-05"
$ws.Range("C10").Value = [double]"1.589797006219889e-123"
$ws.Range("B11").Value = -0.01411993576464326
$ws.Range("C11").Value = 0.1306779273328625
$ws.Range("A12").Value = "DistCenter_pc"
$ws.Range("B12").Value = 0.003132610485043581
$ws.Range("C12").Value = [double]"1.159875405354209e-15"
$ws.Range("B13").Value = 0.4240050243833402
$ws.Range("C13").Value = 0.008487060769297916
$ws.Range("B14").Value = 0.0009047269917179097
$ws.Range("C14").Value = 0.4858056734474153
$ws.Range("B15").Value = -0.002391299024627629
$ws.Range("C15").Value = 0.009970042866253711
$ws.Range("B16").Value = 0.6815411612933205
$ws.Range("C16").Value = [double]"8.062271162743384e-10"
$ws.Range("B17").Value = -0.617924624439722
$ws.Range("C17").Value = [double]"1.807900102246248e-05"
$ws.Range("B18").Value = -0.0001176565290370418
$ws.Range("C18").Value = [double]"2.948447323737908e-23"

$ws = $wb.Worksheets.Item(5)
$ws.Range("B2").Value = -0.8479098035826894
$ws.Range("C2").Value = [double]"2.719245574702528e-06"
$ws.Range("B3").Value = -0.6004508124568247
$ws.Range("C3").Value = [double]"4.979910803543046e-38"
$ws.Range("B4").Value = 0.2649290981427551
$ws.Range("C4").Value = [double]"6.787727647662161e-73"
$ws.Range("B5").Value = 0.0005631358988335549
$ws.Range("C5").Value = 0
$ws.Range("B6").Value = 0.0145505858332803
$ws.Range("C6").Value = [double]"1.194903393616621e-35"
$ws.Range("B7").Value = 0.01653527670301403
$ws.Range("C7").Value = 0.6086000537814176
$ws.Range("B8").Value = 0.4114069667369277
$ws.Range("C8").Value = [double]"3.047826779731331e-18"
$ws.Range("B9").Value = 0.2461229014766139
$ws.Range("C9").Value = [double]"6.565866191086085e-05"
$ws.Range("B10").Value = [double]"-7.136891985176123e-05"
$ws.Range("C10").Value = [double]"3.762361201660085e-115"
$ws.Range("B11").Value = -0.01347359141652811
$ws.Range("C11").Value = 0.1494395297487696
$ws.Range("A12").Value = "DistCenter_pc"
$ws.Range("B12").Value = 0.002864863674671316
$ws.Range("C12").Value = [double]"2.874390645768452e-13"
$ws.Range("B13").Value = 0.3074639580037099
$ws.Range("C13").Value = 0.05627744285071855
$ws.Range("B14").Value = 0.0001972859198704869
$ws.Range("C14").Value = 0.8793325961656545
$ws.Range("B15").Value = -0.003174692911759593
$ws.Range("C15").Value = 0.0006428386938627715
$ws.Range("B16").Value = 0.6395313844347545
$ws.Range("C16").Value = [double]"9.38386687271762e-09"
$ws.Range("B17").Value = -0.7059977701810454
$ws.Range("C17").Value = [double]"1.088753176064531e-06"
$ws.Range("B18").Value = -0.0001274056540285724
$ws.Range("C18").Value = [double]"1.133521154138466e-26"

$ws = $wb.Worksheets.Item(6)
$ws.Range("B2").Value = -0.8399777772111321
$ws.Range("C2").Value = [double]"3.331066930497672e-06"
$ws.Range("B3").Value = -0.6045428841982204
$ws.Range("C3").Value = [double]"1.268259004585874e-38"
$ws.Range("B4").Value = 0.270323617363815
$ws.Range("C4").Value = [double]"4.559004585468925e-75"
$ws.Range("B5").Value = 0.0005528550058776156
$ws.Range("C5").Value = 0
$ws.Range("B6").Value = 0.01395772403569527
$ws.Range("C6").Value = [double]"6.486979691358955e-33"
$ws.Range("B7").Value = 0.01188042354963206
$ws.Range("C7").Value = 0.7124788660224685
$ws.Range("B8").Value = 0.4442370785574933
$ws.Range("C8").Value = [double]"5.791745432460907e-21"
$ws.Range("B9").Value = 0.2648906864191955
$ws.Range("C9").Value = [double]"1.697875304120486e-05"
$ws.Range("B10").Value = [double]"-7.56490167155446e-05"
$ws.Range("C10").Value = [double]"5.349923669723877e-129"
$ws.Range("B11").Value = -0.01591785176280035
$ws.Range("C11").Value = 0.08771295679694188
$ws.Range("A12").Value = "DistCenter_pc"
$ws.Range("B12").Value = 0.003093248777190072
$ws.Range("C12").Value = [double]"3.324265003319426e-15"
$ws.Range("B13").Value = 0.3127625660232409
$ws.Range("C13").Value = 0.05237979571465906
$ws.Range("B14").Value = [double]"3.30781514064416e-05"
$ws.Range("C14").Value = 0.9796588606411736
$ws.Range("B15").Value = -0.003036914322533199
$ws.Range("C15").Value = 0.001047262707183866
$ws.Range("B16").Value = 0.6490487636068379
$ws.Range("C16").Value = [double]"5.371052273617567e-09"
$ws.Range("B17").Value = -0.7866249506715034
$ws.Range("C17").Value = [double]"4.826902590865766e-08"
$ws.Range("B18").Value = -0.0001088860848965803
$ws.Range("C18").Value = [double]"4.052033020188418e-20"

$ws = $wb.Worksheets.Item(7)
$ws.Range("B2").Value = -0.8060665207448
$ws.Range("C2").Value = [double]"7.701025247060786e-06"
$ws.Range("B3").Value = -0.6126103114087215
$ws.Range("C3").Value = [double]"1.117207352013746e-39"
$ws.Range("B4").Value = 0.2560224266135858
$ws.Range("C4").Value = [double]"9.599551121967038e-69"
$ws.Range("B5").Value = 0.0005530239856163506
$ws.Range("C5").Value = 0
$ws.Range("B6").Value = 0.01468889528286265
$ws.Range("C6").Value = [double]"4.244223049928988e-36"
$ws.Range("B7").Value = 0.01110497623448364
$ws.Range("C7").Value = 0.7306446901490857
$ws.Range("B8").Value = 0.4005367305424962
$ws.Range("C8").Value = [double]"3.322567511334661e-17"
$ws.Range("B9").Value = 0.2058693213870702
$ws.Range("C9").Value = 0.0009037176840826875
$ws.Range("B10").Value = [double]"-7.290122764763571e-05"
$ws.Range("C10").Value = [double]"2.575813000244792e-122"
$ws.Range("B11").Value = -0.01805989759582185
$ws.Range("C11").Value = 0.0536019010231162
$ws.Range("A12").Value = "DistCenter_pc"
$ws.Range("B12").Value = 0.003162130698799338
$ws.Range("C12").Value = [double]"9.020536982885302e-16"
$ws.Range("B13").Value = 0.330023249505822
$ws.Range("C13").Value = 0.04126986562042712
$ws.Range("B14").Value = -0.0005189237429694363
$ws.Range("C14").Value = 0.6889159180922912
$ws.Range("B15").Value = -0.00332137826370959
$ws.Range("C15").Value = 0.0003145473298549379
$ws.Range("B16").Value = 0.7631851527544099
$ws.Range("C16").Value = [double]"6.863549269091877e-12"
$ws.Range("B17").Value = -0.6900235554972473
$ws.Range("C17").Value = [double]"1.69850862854673e-06"
$ws.Range("B18").Value = -0.0001191759559142829
$ws.Range("C18").Value = [double]"2.642987032890147e-24"

$ws = $wb.Worksheets.Item(8)
$ws.Range("B2").Value = -0.8901039084776312
$ws.Range("C2").Value = [double]"8.342227721776364e-07"
$ws.Range("B3").Value = -0.5962650984661652
$ws.Range("C3").Value = [double]"2.18499004803133e-37"
$ws.Range("B4").Value = 0.2672607125980139
$ws.Range("C4").Value = [double]"6.998880923574463e-74"
$ws.Range("B5").Value = 0.0005513307950651353
$ws.Range("C5").Value = 0
$ws.Range("B6").Value = 0.01469059266029172
$ws.Range("C6").Value = [double]"3.719189914840156e-36"
$ws.Range("B7").Value = 0.02755426675814544
$ws.Range("C7").Value = 0.3941904312220663
$ws.Range("B8").Value = 0.3954890889703349
$ws.Range("C8").Value = [double]"4.783125189560556e-17"
$ws.Range("B9").Value = 0.2134671759740532
$ws.Range("C9").Value = 0.00051601258597522
$ws.Range("B10").Value = [double]"-7.29366948489388e-05"
$ws.Range("C10").Value = [double]"1.560276284713305e-120"
$ws.Range("B11").Value = -0.004904135117995615
$ws.Range("C11").Value = 0.6028104871244787
$ws.Range("A12").Value = "DistCenter_pc"
$ws.Range("B12").Value = 0.002891614951326176
$ws.Range("C12").Value = [double]"2.896288462037426e-13"
$ws.Range("B13").Value = 0.2765176719126062
$ws.Range("C13").Value = 0.08574553890718964
$ws.Range("B14").Value = 0.0005972246062727779
$ws.Range("C14").Value = 0.6457236344819248
$ws.Range("B15").Value = -0.002812937462276646
$ws.Range("C15").Value = 0.002440694989020094
$ws.Range("B16").Value = 0.6426963435896189
$ws.Range("C16").Value = [double]"7.915728862684651e-09"
$ws.Range("B17").Value = -0.7601018489828213
$ws.Range("C17").Value = [double]"1.348232282191691e-07"
$ws.Range("B18").Value = -0.0001180597180504275
$ws.Range("C18").Value = [double]"2.202841061086458e-23"

$ws = $wb.Worksheets.Item(9)
$ws.Range("B2").Value = -0.9691620992946459
$ws.Range("C2").Value = [double]"8.684254542024193e-08"
$ws.Range("B3").Value = -0.6002884406178015
$ws.Range("C3").Value = [double]"4.960359009080482e-38"
$ws.Range("B4").Value = 0.2684881602706329
$ws.Range("C4").Value = [double]"1.675894374697703e-74"
$ws.Range("B5").Value = 0.0005567191647261914
$ws.Range("C5").Value = 0
$ws.Range("B6").Value = 0.01460714688124717
$ws.Range("C6").Value = [double]"1.273538001538045e-35"
$ws.Range("B7").Value = -0.001910809745769942
$ws.Range("C7").Value = 0.9529037183664893
$ws.Range("B8").Value = 0.4349137635537532
$ws.Range("C8").Value = [double]"4.726778913077635e-20"
$ws.Range("B9").Value = 0.2388011744689241
$ws.Range("C9").Value = 0.0001159568046882343
$ws.Range("B10").Value = [double]"-7.501059248109855e-05"
$ws.Range("C10").Value = [double]"1.764785413028245e-126"
$ws.Range("B11").Value = -0.01715428219868326
$ws.Range("C11").Value = 0.06732419727577339
$ws.Range("A12").Value = "DistCenter_pc"
$ws.Range("B12").Value = 0.003033585615399961
$ws.Range("C12").Value = [double]"1.424784660138011e-14"
$ws.Range("B13").Value = 0.3013213163249902
$ws.Range("C13").Value = 0.06248345445839483
$ws.Range("B14").Value = 0.0007205302232867176
$ws.Range("C14").Value = 0.5783701990343826
$ws.Range("B15").Value = -0.002454218019098535
$ws.Range("C15").Value = 0.008464634534762748
$ws.Range("B16").Value = 0.7183972075373899
$ws.Range("C16").Value = [double]"1.229677497641433e-10"
$ws.Range("B17").Value = -0.7079330061843138
$ws.Range("C17").Value = [double]"1.062129872388704e-06"
$ws.Range("B18").Value = -0.0001199039006267585
$ws.Range("C18").Value = [double]"1.1630012086518e-23"
